$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Move the _GoBack bookmark from the first paragraph ("To: File") to
#    the end of the "Backing: ..." paragraph (after "... affirmed by
#    the members.").
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$backingPara = $d.Paragraphs(10)
# Sanity check we have the right paragraph before mutating it.
if ($backingPara.Range.Text -notmatch "^Backing:") {
    throw "Unexpected paragraph 10 content: $($backingPara.Range.Text)"
}

$marker = "zzGoBackMarkerzz"
$backingPara.Range.InsertAfter($marker)

$markerRange = $d.Content
$markerRange.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markerRange)
$d.Bookmarks("_GoBack").Range.Text = ""

# ---------------------------------------------------------------------
# 2. "Claim: Every member ..." -- insert "that is " before "derived
#    from their".
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Claim: Every member of an organization has a warrant to influence organizational objectives and effectuate organizational change to achieve those objectives derived from their",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Claim: Every member of an organization has a warrant to influence organizational objectives and effectuate organizational change to achieve those objectives that is derived from their",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. "The observations of the theorists ..." -- reword / restructure.
#    Done as two Find/Replace calls so the untouched, proofErr-wrapped
#    "according to" run in between is left alone.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The observations of the theorists in the human relations movement, such as the Hawthorne studies, demonstrate that organization members do not absentmindedly react to stimulus in the organizational environment but respond ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The observations of the theorists in the human relations movement demonstrate that organization members do not absentmindedly react to stimulus in the environment but respond ",
    2) | Out-Null

$d.Content.Find.Execute(
    " the meaning they assign to events in the organizational environment (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " the meaning they assign to events that occur in the organizational setting (",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4. "The machine metaphor of organization ..." -- reword.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "The machine metaphor of organization is little more than indentured servitude, which does not conform to the values of American society.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The machine metaphor of organization is little more than indentured servitude, which is inconsistent with the ideals and values of American society.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5. "Relying solely on fiduciary leaders ..." -- reword.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Relying solely on fiduciary leaders to effectuate organization change increases the risk of catastrophic failure because fiduciary leaders are fallible humans constrained by bounded rationality.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Relying solely on fiduciary leaders to effectuate change in an organization increases the risk of catastrophic failure because fiduciary leaders are fallible humans constrained by bounded rationality.",
    2) | Out-Null

Write-Output "Done"
